$wb = $excel.ActiveWorkbook

# Sheet1: remove the 3 data rows (Nguyen Van Dinh1/2/3), keeping only the
# header row. This also drops the now-unused shared strings.
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Rows("2:4").Delete()

# Sheet2: add a new value 5 in B2 (C1 keeps its "tuoi" text; its shared
# string index is simply renumbered once the unused strings are removed).
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("B2").Value = 5
